$d = $word.ActiveDocument

$replacements = @(
    @{old="319÷6="; new="716÷2="},
    @{old="233÷2="; new="284÷9="},
    @{old="883÷3="; new="844÷7="},
    @{old="571÷8="; new="642÷5="},
    @{old="163÷7="; new="887÷6="},
    @{old="223÷7="; new="888÷2="},
    @{old="300÷3="; new="409÷3="},
    @{old="511÷7="; new="961÷7="},
    @{old="220÷8="; new="733÷5="},
    @{old="347÷7="; new="588÷2="},
    @{old="855÷7="; new="306÷9="},
    @{old="471÷2="; new="358÷9="},
    @{old="195÷4="; new="181÷8="},
    @{old="804÷6="; new="607÷8="},
    @{old="105÷2="; new="489÷2="},
    @{old="670÷7="; new="787÷6="},
    @{old="737÷3="; new="648÷5="},
    @{old="151÷4="; new="251÷4="},
    @{old="750÷5="; new="547÷2="},
    @{old="417÷6="; new="410÷3="},
    @{old="349÷2="; new="412÷5="},
    @{old="567÷8="; new="398÷7="},
    @{old="400÷9="; new="231÷3="},
    @{old="111÷8="; new="995÷5="},
    @{old="233÷4="; new="755÷4="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
